$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.230.37"
$ws.Range("E2").Value = "  +0.52%  "
$ws.Range("D3").Value = "1.844.07"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'244.43"
$ws.Range("E5").Value = "  -0.06%  "
$ws.Range("D6").Value = "'0.6290"
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").Value = "'0.07548"
$ws.Range("E8").Value = "  -0.33%  "
$ws.Range("D9").Value = "'0.2954"
$ws.Range("E9").Value = "  +0.12%  "
$ws.Range("E10").Value = "  +2.12%  "
$ws.Range("D11").Value = "'0.07724"
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("D12").Value = "1.843.34"
$ws.Range("E12").Value = "  +0.40%  "
$ws.Range("D13").Value = "'5.037"
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "'0.6802"
$ws.Range("E14").Value = "  +1.19%  "
$ws.Range("D15").Value = "'83.46"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'0.000009306"
$ws.Range("E16").Value = "  -3.08%  "
$ws.Range("D17").Value = "'5.995"
$ws.Range("E17").Value = "  -2.14%  "
$ws.Range("D18").Value = "29.216.26"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "2.094.12"
$ws.Range("E19").Value = "  -0.35%  "
$ws.Range("D20").Value = "'234.52"
$ws.Range("E20").Value = "  +3.32%  "
$ws.Range("D21").Value = "'12.75"
$ws.Range("E21").Value = "  +1.15%  "
$ws.Range("D22").Value = "'1.004"
$ws.Range("E22").Value = "  +0.31%  "
$ws.Range("D23").Value = "'7.193"
$ws.Range("E23").Value = "  -0.51%  "
$ws.Range("D24").Value = "'1.004"
$ws.Range("E24").Value = "  +0.24%  "
$ws.Range("D25").Value = "'160.81"
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "'0.1404"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").Value = "'8.577"
$ws.Range("E27").Value = "  +0.23%  "
$ws.Range("D28").Value = "'18.02"
$ws.Range("D29").Value = "'1.501"
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +1.75%  "
$ws.Range("D31").Value = "'4.165"
$ws.Range("E31").Value = "  +1.86%  "
$ws.Range("D32").Value = "'0.05593"
$ws.Range("E32").Value = "  +3.02%  "
$ws.Range("D33").Value = "'1.210"
$ws.Range("E33").Value = "  +0.37%  "
$ws.Range("D34").Value = "'0.7539"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "'1.863"
$ws.Range("E35").Value = "  -0.19%  "
$ws.Range("D36").Value = "'1.153"
$ws.Range("E36").Value = "  +0.78%  "
$ws.Range("D37").Value = "'2.666"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "1.241.73"
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'2.774"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'0.01798"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'6.651"
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("D42").Value = "'0.9058"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("D43").Value = "'1.003"
$ws.Range("E43").Value = "  +0.08%  "
$ws.Range("D44").Value = "'102.71"
$ws.Range("E44").Value = "  +0.36%  "
$ws.Range("D45").Value = "1.996.07"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("D46").Value = "'66.81"
$ws.Range("E46").Value = "  +2.54%  "
$ws.Range("D47").Value = "'0.5103"
$ws.Range("E47").Value = "  -0.27%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("D49").Value = "'0.4112"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "'9.160"
$ws.Range("E50").Value = "  +0.09%  "
$ws.Range("D51").Value = "'0.05844"
$ws.Range("E51").Value = "  +1.08%  "
